$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.501534223556519
$ws.Range("B1").Value = 2.264607429504395
$ws.Range("C1").Value = 1.99350106716156
$ws.Range("D1").Value = 1.74684751033783
$ws.Range("E1").Value = 1.282963156700134
